# CoreAcMain.xlsx - DBD sheet edits
# - G10 gains a third line "空白三位:合計" appended to the existing
#   "00A:傳統帳冊 / 201:利變年金帳冊" note.
# - G11 (previously empty) gets a new two-line note "TOL:合計 / NTD:新台幣".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")

$ws.Range("G10").Value = "00A:傳統帳冊`n201:利變年金帳冊`n空白三位:合計"
$ws.Range("G11").Value = "TOL:合計`nNTD:新台幣"

# Row heights grow to fit the (now) 3-line / 2-line wrapped text.
$ws.Rows.Item(10).RowHeight = 48.6
$ws.Rows.Item(11).RowHeight = 32.4

# Leave the selection where the edit actually happened.
$ws.Range("G11").Select()
